# Insert a new row above row 6 to hold "Bonus Languages" data for each
# Divine Domain, defaulting to "None" (archetypes may now grant bonus
# languages instead).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X")
foreach ($col in $cols) {
    $ws.Range("$col`6").Value = "None"
}
